$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# feat: enchant world goldAdd buff
#
# The "EnchantType" legend cell (J3) gains three new categories (world 1/2/3
# gold bonuses), and the newly-added enchant rows (5-19, which model the
# three world-gold-bonus entries at levels 1-5) get re-pointed from the old
# placeholder EnchantType code (2 = "金币加成") onto their own new codes:
#   rows 5-9   -> 9  (第一世界金币加成 / world-1 gold bonus)
#   rows 10-14 -> 10 (第二世界金币加成 / world-2 gold bonus)
#   rows 15-19 -> 11 (第三世界金币加成 / world-3 gold bonus)
# ---------------------------------------------------------------------------

$j3Text = @"
0.其他
1.伤害加成
2.金币加成
3.钻石加成
4.暴击加成
5.移动速度加成
6.宝箱伤害加成
7.四级金币宝箱加成
8.倍率资源加成金币
9.第一世界金币加成
10.第二世界金币加成
11.第三世界金币加成
"@
$ws.Range("J3").Value = $j3Text

# Row 3 has to grow taller to fit the 3 extra legend lines.
$ws.Rows.Item(3).RowHeight = 255

# Column J (EnchantType) is widened a touch to fit the new two-digit codes.
$ws.Columns.Item(10).ColumnWidth = 15.5

# Re-point the EnchantType column for the new goldAdd-buff rows.
$ws.Range("J5:J9").Value = 9
$ws.Range("J10:J14").Value = 10
$ws.Range("J15:J19").Value = 11

# Leave the selection on the cell that was edited.
$ws.Range("J3").Select()
